$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column K (complete) so the existing
# "complete" and "lastedit" columns shift right by one, then populate
# the new column K with the "people_rgb" header/value.
$ws.Range("K1").EntireColumn.Insert()

$ws.Range("K1").Value = "people_rgb"

# Update people_xy coordinates for row 2
$ws.Range("J2").Value = "(176,94)|(478,589)|(783,136)|"

# Populate the new people_rgb column for row 2
$ws.Range("K2").Value = "(0, 0, 0, 1)|(0, 0, 0, 1)|(0, 0, 0, 1)|"

# Update the lastedit timestamp (now in column M after the insert)
$ws.Range("M2").Value = "05/01/2023 16:30:13"
